# Generate Report for Handback
#
# - Marks handoff rows as "Handed back: in sync with en-US" (was "Ready for
#   handoff") on the Overview sheet and on each language sheet.
# - Stamps the "Latest Handback DateTime" (column H) for each language sheet.
# - Fills in "Latest Target File" (F) / "Latest Handback File" (G) columns
#   with the handed-back file names + hyperlinks, for each data row.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: Status columns (B = zh-cn, C = de-de) for each file row.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusText
$overview.Range("C2").Value = $statusText
$overview.Range("B3").Value = $statusText
$overview.Range("C3").Value = $statusText

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# Status column
$zh.Range("C2").Value = $statusText
$zh.Range("C3").Value = $statusText

# Latest Handback DateTime
$zh.Range("H2").Value = "2016-03-12 00:32:30"
$zh.Range("H3").Value = "2016-03-12 00:32:30"

# Latest Target File / Latest Handback File (+ hyperlinks)
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/ca2b86beeb70ed70127b5a77b0dbf43bd73b4a72/e2e/75c0ea31-9b7c-49de-8961-cfe8e4c994ff.md", "", "", "75c0ea31-9b7c-49de-8961-cfe8e4c994ff.md")
$zh.Hyperlinks.Add($zh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/39a6c0f548dcb534b2890d01a48bf7b100f34378/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/hb/75c0ea31-9b7c-49de-8961-cfe8e4c994ff.1f8af19720502e48cf8fb9bd86ccbcaa9199fde9.zh-cn.xlf", "", "", "75c0ea31-9b7c-49de-8961-cfe8e4c994ff.1f8af19720502e48cf8fb9bd86ccbcaa9199fde9.zh-cn.xlf")

$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/ca2b86beeb70ed70127b5a77b0dbf43bd73b4a72/e2e/f86032e9-c9e9-421d-a86d-6d774bc12014.md", "", "", "f86032e9-c9e9-421d-a86d-6d774bc12014.md")
$zh.Hyperlinks.Add($zh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/39a6c0f548dcb534b2890d01a48bf7b100f34378/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/hb/f86032e9-c9e9-421d-a86d-6d774bc12014.b73e53f97c6f07904a8bb7469bf17a0de1c5a033.zh-cn.xlf", "", "", "f86032e9-c9e9-421d-a86d-6d774bc12014.b73e53f97c6f07904a8bb7469bf17a0de1c5a033.zh-cn.xlf")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

# Status column
$de.Range("C2").Value = $statusText
$de.Range("C3").Value = $statusText

# Latest Handback DateTime
$de.Range("H2").Value = "2016-03-12 00:32:35"
$de.Range("H3").Value = "2016-03-12 00:32:35"

# Latest Target File / Latest Handback File (+ hyperlinks)
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/ca2b86beeb70ed70127b5a77b0dbf43bd73b4a72/e2e/75c0ea31-9b7c-49de-8961-cfe8e4c994ff.md", "", "", "75c0ea31-9b7c-49de-8961-cfe8e4c994ff.md")
$de.Hyperlinks.Add($de.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/568868c54c194088a248aa68ae5d89c2272bb2d4/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/hb/75c0ea31-9b7c-49de-8961-cfe8e4c994ff.1f8af19720502e48cf8fb9bd86ccbcaa9199fde9.de-de.xlf", "", "", "75c0ea31-9b7c-49de-8961-cfe8e4c994ff.1f8af19720502e48cf8fb9bd86ccbcaa9199fde9.de-de.xlf")

$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/ca2b86beeb70ed70127b5a77b0dbf43bd73b4a72/e2e/f86032e9-c9e9-421d-a86d-6d774bc12014.md", "", "", "f86032e9-c9e9-421d-a86d-6d774bc12014.md")
$de.Hyperlinks.Add($de.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/568868c54c194088a248aa68ae5d89c2272bb2d4/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/hb/f86032e9-c9e9-421d-a86d-6d774bc12014.b73e53f97c6f07904a8bb7469bf17a0de1c5a033.de-de.xlf", "", "", "f86032e9-c9e9-421d-a86d-6d774bc12014.b73e53f97c6f07904a8bb7469bf17a0de1c5a033.de-de.xlf")
